# Active_Outages.xlsx - add new outage row to the "R1" sheet (first sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank-but-present text cells: Excel stores these as an explicit empty
# string (t="str" / <v/>), not as an absent cell. Assigning "" directly
# clears the cell entirely, so force a text entry the way a user would
# (leading apostrophe = "treat as text"), then strip the resulting
# quote-prefix style so the cell formatting matches the rest of the sheet.
$blankCells = @("A5", "C5", "E5", "F5", "G5", "H5", "K5")
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("B5").Value = "R4"
$ws.Range("D5").Value = "JED0123"
$ws.Range("I5").Value = "SCECO"
$ws.Range("J5").Value = "In progress"
$ws.Range("L5").Value = "Latis"
